$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 ("Task 13: Polish form sizes for different screens"):
#   Status moves from Pending -> In Development (copy style from the
#   "In Development" legend cell G3, then set the text).
$ws.Range("G3").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = "In Development"

# New "Projects" CRUD tasks - fill in the full names for rows 14-18
# (a project-owner column is being added to the Projects table, so the
# todo list gets the matching Create/Read/Update/Delete task rows).
$ws.Range("A14").Value = "Task 14: Projects Page Frontend"
$ws.Range("A15").Value = "Task 15: Projects Cread (CRUD)"
$ws.Range("A16").Value = "Task 16: Projects Read (CRUD)"
$ws.Range("A17").Value = "Task 17: Projects Update (CRUD)"
$ws.Range("A18").Value = "Task 18: Projects Delete (CRUD)"

# Row 18 gets a Status of "Pending" too, matching rows 14-17's style
$ws.Range("B17").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "Pending"

# Record a new "Date Last Updated" value for row 13 in column D.
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "02/22/2024"
$ws.Range("D13").ClearFormats()

$ws.Range("C11").Select()
